$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text values so Excel does not reinterpret strings like "27.510.90"
# or "  -0.13%  " as numbers/dates/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.510.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.84%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.613.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.61%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.515.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.46"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.88%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.96"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.38%  "

# Rows 26/27 swap Cosmos <-> Stellar (with updated figures)
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.111"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.83"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.443.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.24%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.02%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.24"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.72%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.37%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.758.05"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.23"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.35%  "
